# "Updated cryptos list" refresh: bumps the Price (column D) and
# Volume(1h) (column E) figures for each coin row, and swaps the
# ARBITRUM / ImmutableX rows (32 and 33) to their new order/prices.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D retains text formatting so numeric-looking values
# (e.g. "315.53") are not auto-converted to numbers, matching the
# original inline-string cell type.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '27.379.85'
$ws.Range('E2').Value = '  -1.51%  '
$ws.Range('D3').Value = '1.827.87'
$ws.Range('E3').Value = '  -1.69%  '
$ws.Range('D4').Value = '1.003'
$ws.Range('E4').Value = '  -3.32%  '
$ws.Range('D5').Value = '315.53'
$ws.Range('E5').Value = '  -2.54%  '
$ws.Range('E6').Value = '  -2.93%  '
$ws.Range('D7').Value = '0.4298'
$ws.Range('E7').Value = '  -2.61%  '
$ws.Range('D8').Value = '0.3701'
$ws.Range('E8').Value = '  -3.05%  '
$ws.Range('D9').Value = '0.07255'
$ws.Range('E9').Value = '  -2.62%  '
$ws.Range('D10').Value = '0.8668'
$ws.Range('E10').Value = '  -2.44%  '
$ws.Range('D11').Value = '21.16'
$ws.Range('E11').Value = '  -2.07%  '
$ws.Range('D12').Value = '1.819.42'
$ws.Range('E12').Value = '  -2.21%  '
$ws.Range('D13').Value = '6.679'
$ws.Range('E13').Value = '  -1.01%  '
$ws.Range('D14').Value = '5.352'
$ws.Range('E14').Value = '  -3.44%  '
$ws.Range('D15').Value = '0.07070'
$ws.Range('E15').Value = '  -2.01%  '
$ws.Range('D16').Value = '87.94'
$ws.Range('E16').Value = '  +1.99%  '
$ws.Range('E17').Value = '  -3.29%  '
$ws.Range('D18').Value = '0.000008898'
$ws.Range('E18').Value = '  -2.49%  '
$ws.Range('D19').Value = '1.003'
$ws.Range('E19').Value = '  -3.00%  '
$ws.Range('D20').Value = '15.26'
$ws.Range('E20').Value = '  -2.17%  '
$ws.Range('D21').Value = '27.390.41'
$ws.Range('E21').Value = '  -1.52%  '
$ws.Range('D22').Value = '5.165'
$ws.Range('E22').Value = '  -2.75%  '
$ws.Range('D23').Value = '10.87'
$ws.Range('E23').Value = '  -3.69%  '
$ws.Range('D24').Value = '2.056.91'
$ws.Range('E24').Value = '  -1.56%  '
$ws.Range('E25').Value = '  -2.67%  '
$ws.Range('D26').Value = '153.40'
$ws.Range('E26').Value = '  -3.77%  '
$ws.Range('D27').Value = '18.42'
$ws.Range('E27').Value = '  -1.95%  '
$ws.Range('D28').Value = '2.139'
$ws.Range('E28').Value = '  +6.82%  '
$ws.Range('D29').Value = '5.290'
$ws.Range('E29').Value = '  -1.76%  '
$ws.Range('D30').Value = '117.26'
$ws.Range('E30').Value = '  -1.28%  '
$ws.Range('D31').Value = '0.08835'
$ws.Range('E31').Value = '  -3.08%  '
$ws.Range('B32').Value = 'ARBITRUM'
$ws.Range('C32').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D32').Value = '1.204'
$ws.Range('E32').Value = '  -1.23%  '
$ws.Range('B33').Value = 'ImmutableX'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D33').Value = '0.7667'
$ws.Range('E33').Value = '  -1.31%  '
$ws.Range('D34').Value = '4.507'
$ws.Range('E34').Value = '  -2.19%  '
$ws.Range('D35').Value = '2.865'
$ws.Range('E35').Value = '  -5.27%  '
$ws.Range('E36').Value = '  -3.16%  '
$ws.Range('E37').Value = '  -3.28%  '
$ws.Range('D38').Value = '0.01961'
$ws.Range('E38').Value = '  -1.38%  '
$ws.Range('D39').Value = '0.05283'
$ws.Range('E39').Value = '  -0.80%  '
$ws.Range('D40').Value = '2.882'
$ws.Range('E40').Value = '  +0.35%  '
$ws.Range('D41').Value = '7.146'
$ws.Range('E41').Value = '  +2.51%  '
$ws.Range('E42').Value = '  -0.18%  '
$ws.Range('D43').Value = '0.5073'
$ws.Range('E43').Value = '  -2.86%  '
$ws.Range('D44').Value = '8.667'
$ws.Range('E44').Value = '  -1.46%  '
$ws.Range('D45').Value = '10.59'
$ws.Range('E45').Value = '  -1.84%  '
$ws.Range('D46').Value = '106.35'
$ws.Range('E46').Value = '  -4.27%  '
$ws.Range('D47').Value = '0.4730'
$ws.Range('E47').Value = '  -0.19%  '
$ws.Range('D48').Value = '0.06421'
$ws.Range('E48').Value = '  -2.47%  '
$ws.Range('E49').Value = '  -3.27%  '
$ws.Range('E50').Value = '  -2.86%  '
$ws.Range('D51').Value = '1.825'
$ws.Range('E51').Value = '  -3.06%  '
